$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E, shifting Sequence.R..R.Stop (old E:H) to F:I
$ws.Columns("E:E").Insert()

# Give the new column roughly the same width as column D
$ws.Range("E1").ColumnWidth = 28.5

# Header for the new "Group" column
$ws.Range("E1").Value = "Group"

# Fill the new "Group" column with the taxonomic group, copying the format
# of the adjacent Sequence.R cell on each row
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = "Echinodermata"

$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "Echinodermata"

# Move the active selection to A2
$ws.Range("A2").Select()
